# msz - field hint and error checks part 3 + listener experiments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new test-case row (row 6) mirroring the existing table layout.
$ws.Range("A6").Value = "102_AutomobileInsurance_003_InsurantData_002_FieldHintsAndErrors"
$ws.Range("B6").Value = "Button Next from Page VehicleData"
$ws.Range("C6").Value = "Insurant Page check for hints regarding mandatory fields"
$ws.Range("D6").Value = "102_AutomobileInsurance_003_InsurantData_002_EnterValuesInWrongFormat"
$ws.Range("E6").Value = "Insurant Page check error hint formatting"
$ws.Range("F6").Value = "102_AutomobileInsurance_003_InsurantData_002_EnterValuesInWrongFormat Part 2"
$ws.Range("G6").Value = "Insurant Page check error hint formatting Part 2"

# Column F now holds a longer string than before, so widen it to fit (best-fit).
$ws.Columns.Item(6).AutoFit()

# Move the active selection, as the author left off around F14.
$ws.Range("F14").Select()
